# Lab Exercise - Docker Fundamental Commands
# The author trimmed the document so that it stops right after the
# "Step 6: Remove Images" / "Note: ..." section.  Concretely:
#   - The paragraph that only holds a manual page-break run
#     (<w:br w:type="page"/>) is kept, but the run itself is removed,
#     leaving an empty paragraph (same as the very last paragraph of
#     the trimmed document elsewhere in the file).
#   - Every paragraph after that page-break paragraph ("3. Working
#     with Custom Images" through the final empty "Running a simple
#     web app in Docker" list filler paragraph) is deleted outright.
#   - The trailing <w:sectPr> (page setup) is left untouched.

$d = $word.ActiveDocument

# Locate the heading that starts the section to be removed ("3. Working
# with Custom Images"). Using Find (rather than a hard-coded paragraph
# index) keeps this robust to any earlier, unrelated edits.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "3. Working with Custom Images"
$find.Forward = $true
$find.Wrap = 1
$null = $find.Execute()

$headingPara = $d.Range($find.Parent.Start, $find.Parent.Start).Paragraphs(1)

# 1) Delete everything from the start of that heading paragraph through
#    the end of the document's main story. That removes every paragraph
#    from "3. Working with Custom Images" down through the last (empty)
#    paragraph of the "Lab Summary" bullet list, while leaving the
#    trailing <w:sectPr> (section/page setup) untouched.
$tail = $d.Range($headingPara.Range.Start, $d.Content.End)
$tail.Delete()

# 2) The paragraph immediately preceding the removed section holds only
#    a manual page break (<w:br w:type="page"/>). Strip that run out but
#    keep the paragraph mark, leaving an empty paragraph in its place.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$runRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
if ($runRange.Text.Length -gt 0) {
    $runRange.Delete()
}
